$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.268.82'
$ws.Range("E2").Value = '  -2.71%  '

$ws.Range("D3").Value = '1.565.08'
$ws.Range("E3").Value = '  -3.58%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.39'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.476'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.43%  '

$ws.Range("E8").Value = '  -2.55%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0606'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.11%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.79'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.73%  '

$ws.Range("E11").Value = '  -0.78%  '

$ws.Range("D12").Value = '1.781.41'
$ws.Range("E12").Value = '  -3.62%  '

$ws.Range("D13").Value = '1.565.43'
$ws.Range("E13").Value = '  -3.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.53%  '

$ws.Range("E15").Value = '  -3.34%  '

$ws.Range("D16").Value = '25.267.11'
$ws.Range("E16").Value = '  -2.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '59.36'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.00%  '

$ws.Range("D18").Value = '0.0₃0713'
$ws.Range("E18").Value = '  -3.44%  '

$ws.Range("E19").Value = '  -0.40%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '185.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.10%  '

$ws.Range("E21").Value = '  -2.60%  '

$ws.Range("E22").Value = '  -2.91%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.87'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.22%  '

$ws.Range("E24").Value = '  -2.26%  '

$ws.Range("E25").Value = '  -0.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '139.76'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.69%  '

$ws.Range("E27").Value = '  -7.60%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '14.86'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.00%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.45'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.83%  '

$ws.Range("E30").Value = '  -6.27%  '

$ws.Range("E31").Value = '  -3.61%  '

$ws.Range("E32").Value = '  -2.94%  '

$ws.Range("E33").Value = '  -3.83%  '

$ws.Range("E34").Value = '  -1.65%  '

$ws.Range("E35").Value = '  -4.18%  '

$ws.Range("D36").Value = '1.086.19'
$ws.Range("E36").Value = '  -3.53%  '

$ws.Range("E37").Value = '  -0.82%  '

$ws.Range("E38").Value = '  -4.68%  '

$ws.Range("E39").Value = '  -2.16%  '

$ws.Range("E40").Value = '  -3.83%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.771'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.68%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.761'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.32%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '93.28'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.44%  '

$ws.Range("E44").Value = '  -2.49%  '

$ws.Range("D45").Value = '1.695.87'
$ws.Range("E45").Value = '  -3.61%  '

$ws.Range("D46").Value = '0.0₆0111'
$ws.Range("E46").Value = '  -2.60%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '52.70'
$ws.Range("D47").Style = "Normal"

$ws.Range("E48").Value = '  -3.65%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.42'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.79%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.407'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.73%  '

$ws.Range("E51").Value = '  -0.52%  '
